$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert a new row above current row 2, shifting existing data rows down.
$ws.Rows.Item(2).Insert()

# The inserted row picks up the bold/bordered header formatting; clear it
# back to the plain/default style used by the other data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new student's data.
$ws.Cells.Item(2, 1).Value = 18330051920026
$ws.Cells.Item(2, 2).Value = "ROJAS"
$ws.Cells.Item(2, 3).Value = "MAZA"
$ws.Cells.Item(2, 4).Value = "ANGEL GABRIEL"
$ws.Cells.Item(2, 5).Value = "MATEMÁTICAS APLICADAS"
$ws.Cells.Item(2, 6).Value = "6AEV"
$ws.Cells.Item(2, 7).Value = 2
